# Apply edits described by the diff: add a new date column (AI on "data", AH on "pocetR")
# with updated 12. 10. 2021 figures, and bump the "aktualizace" date in the footer rows.
$wb = $excel.ActiveWorkbook

# --- Sheet "data": new column AI ---
$ws1 = $wb.Worksheets.Item("data")

# Header cell AI1: copy the format from AH1 (which already has the header style) then set the date text.
$ws1.Range("AH1").Copy($ws1.Range("AI1"))
$ws1.Range("AI1").Value = "12. 10. 2021"

# Data rows 2-45: new numeric values for 12. 10. 2021
$ws1.Cells.Item(2, 35).Value = 0.2
$ws1.Cells.Item(3, 35).Value = 0.12
$ws1.Cells.Item(4, 35).Value = 0.41
$ws1.Cells.Item(5, 35).Value = 0.25
$ws1.Cells.Item(6, 35).Value = 0.16
$ws1.Cells.Item(7, 35).Value = 0.19
$ws1.Cells.Item(8, 35).Value = 0.21
$ws1.Cells.Item(9, 35).Value = 0.2
$ws1.Cells.Item(10, 35).Value = 0.21
$ws1.Cells.Item(11, 35).Value = 0.19
$ws1.Cells.Item(12, 35).Value = 0.2
$ws1.Cells.Item(13, 35).Value = 0.32
$ws1.Cells.Item(14, 35).Value = 0.19
$ws1.Cells.Item(15, 35).Value = 0.18
$ws1.Cells.Item(16, 35).Value = 0.23
$ws1.Cells.Item(17, 35).Value = 0.2
$ws1.Cells.Item(18, 35).Value = 0.24
$ws1.Cells.Item(19, 35).Value = 0.24
$ws1.Cells.Item(20, 35).Value = 0.19
$ws1.Cells.Item(21, 35).Value = 0.13
$ws1.Cells.Item(22, 35).Value = 0.12
$ws1.Cells.Item(23, 35).Value = 0.22
$ws1.Cells.Item(24, 35).Value = 0.41
$ws1.Cells.Item(25, 35).Value = 0.39
$ws1.Cells.Item(26, 35).Value = 0.11
$ws1.Cells.Item(27, 35).Value = 0.09
$ws1.Cells.Item(28, 35).Value = 0.14
$ws1.Cells.Item(29, 35).Value = 0.22
$ws1.Cells.Item(30, 35).Value = 0.08
$ws1.Cells.Item(31, 35).Value = 0.11
$ws1.Cells.Item(32, 35).Value = 0.13
$ws1.Cells.Item(33, 35).Value = 0.22
$ws1.Cells.Item(34, 35).Value = 0.18
$ws1.Cells.Item(35, 35).Value = 0.08
$ws1.Cells.Item(36, 35).Value = 0.13
$ws1.Cells.Item(37, 35).Value = 0.13
$ws1.Cells.Item(38, 35).Value = 0.08
$ws1.Cells.Item(39, 35).Value = 0.27
$ws1.Cells.Item(40, 35).Value = 0.13
$ws1.Cells.Item(41, 35).Value = 0.08
$ws1.Cells.Item(42, 35).Value = 0.08
$ws1.Cells.Item(43, 35).Value = 0.09
$ws1.Cells.Item(44, 35).Value = 0.18
$ws1.Cells.Item(45, 35).Value = 0.31

# Row 46: bump the "aktualizace" date in the footer label
$ws1.Range("A46").Value = "Život během pandemie, Strategie domácností, % respondentů celkově a ve skupinách, aktualizace 20. 10. 2021"

# --- Sheet "pocetR": new column AH ---
$ws2 = $wb.Worksheets.Item("pocetR")

# Header cell AH1: copy the format from AG1 then set the date text.
$ws2.Range("AG1").Copy($ws2.Range("AH1"))
$ws2.Range("AH1").Value = "12. 10. 2021"

# Data rows 2-23: new numeric values for 12. 10. 2021
$ws2.Cells.Item(2, 34).Value = 1836
$ws2.Cells.Item(3, 34).Value = 187
$ws2.Cells.Item(4, 34).Value = 352
$ws2.Cells.Item(5, 34).Value = 1297
$ws2.Cells.Item(6, 34).Value = 887
$ws2.Cells.Item(7, 34).Value = 163
$ws2.Cells.Item(8, 34).Value = 525
$ws2.Cells.Item(9, 34).Value = 261
$ws2.Cells.Item(10, 34).Value = 850
$ws2.Cells.Item(11, 34).Value = 153
$ws2.Cells.Item(12, 34).Value = 111
$ws2.Cells.Item(13, 34).Value = 722
$ws2.Cells.Item(14, 34).Value = 850
$ws2.Cells.Item(15, 34).Value = 625
$ws2.Cells.Item(16, 34).Value = 361
$ws2.Cells.Item(17, 34).Value = 193
$ws2.Cells.Item(18, 34).Value = 666
$ws2.Cells.Item(19, 34).Value = 611
$ws2.Cells.Item(20, 34).Value = 245
$ws2.Cells.Item(21, 34).Value = 561
$ws2.Cells.Item(22, 34).Value = 323
$ws2.Cells.Item(23, 34).Value = 157

# Row 24: updated label, plus an empty placeholder cell in the new column
# (keeps the blank-row pattern used by the rest of the footer row: B24:AG24 are blank cells)
$ws2.Range("A24").Value = "Život během pandemie, Strategie domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 20. 10. 2021"
$ws2.Range("AG24").Copy($ws2.Range("AH24"))

